$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the existing "adj" keyword row as Enabled = TRUE
$ws.Range("B2").Value = $true

# Add a new keyword row: "Bahi", Enabled = FALSE
$ws.Range("A3").Value = "Bahi"
$ws.Range("B3").Value = $false

# Move the active selection back to A1 (matches the saved view state)
$ws.Range("A1").Select()
